$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Data modifikation" (sheet3): add the example "rename variable" row
# showing a long tracking URL (G3, hyperlinked) being renamed to a shorter
# value (H3, plain text).
# ---------------------------------------------------------------------------
$wsMod = $wb.Worksheets.Item("Data modifikation")

$longUrl  = "https://trafikkort.vejdirektoratet.dk/?utm_medium=newsletter_ubivox&utm_source=20180308_Sibirisk%20vinter%20skabte%20travlhed%20for%20sneryddere%20og%20saltspredere&utm_campaign=Sibirisk%20vinter%20skabte%20travlhed%20for%20sneryddere%20og%20saltspredere"
$shortUrl = "https://trafikkort.vejdirektoratet.dk/?utm_medium=newsletter_ubivox&utm_source=20180308_Sibirisk%20vinter"

$wsMod.Range("G3").Value = $longUrl
$wsMod.Range("H3").Value = $shortUrl
$wsMod.Hyperlinks.Add($wsMod.Range("G3"), $longUrl) | Out-Null

$wsMod.Range("H3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Data modifikation 2" (sheet4): extend the variable start/end mapping
# table with an exclude-list column and several new example rows.
# ---------------------------------------------------------------------------
$wsMod2 = $wb.Worksheets.Item("Data modifikation 2")

# Match column B's width to column A, and give the new column C a width.
$colAWidth = $wsMod2.Range("A1").EntireColumn.ColumnWidth
$wsMod2.Range("B1").EntireColumn.ColumnWidth = $colAWidth
$wsMod2.Range("C1").EntireColumn.ColumnWidth = 46.33

# Header row - new column C header.
$wsMod2.Range("C2").Value = "Ekskluder værdier der slutter med: (Separer med ;)"

$excludeList = "GroupA,GroupG;GroupA1"
$oldUrl      = "https://trafikkort.vejdirektoratet.dk/?utm_source=danwest.de&utm_campaign=3ec1f9fb03-EMAIL_CAMPAIGN_2018_06_11_08"

$urlVisibleGroups      = "https://trafikkort.vejdirektoratet.dk/index.html?visiblegroups"
$urlUserType2          = "https://trafikkort.vejdirektoratet.dk/index.html?usertype=2"
$urlUserType3          = "https://trafikkort.vejdirektoratet.dk/index.html?usertype=3"
$urlVisibleGroupsA1    = "https://trafikkort.vejdirektoratet.dk/index.html?visibleGroups=GroupA1"
$urlShow               = "https://trafikkort.vejdirektoratet.dk/index.html?show"
$urlLat                = "https://trafikkort.vejdirektoratet.dk/index.html?lat"
$urlGclid              = "https://trafikkort.vejdirektoratet.dk/?gclid"

# Row 3 - new example entry, both columns hyperlinked.
$wsMod2.Range("A3").Value = $urlVisibleGroups
$wsMod2.Range("B3").Value = $urlVisibleGroups
$wsMod2.Range("C3").Value = $excludeList
$wsMod2.Hyperlinks.Add($wsMod2.Range("A3"), $urlVisibleGroups) | Out-Null
$wsMod2.Hyperlinks.Add($wsMod2.Range("B3"), $urlVisibleGroups) | Out-Null

# Row 4 - new example entry.
$wsMod2.Range("A4").Value = $urlUserType2
$wsMod2.Range("B4").Value = $urlUserType2

# Row 5 - the pre-existing example entry, shifted down, with exclude list added.
$wsMod2.Range("A5").Value = $oldUrl
$wsMod2.Range("B5").Value = $oldUrl
$wsMod2.Range("C5").Value = $excludeList

# Row 6 - new example entry.
$wsMod2.Range("A6").Value = $urlUserType3
$wsMod2.Range("B6").Value = $urlUserType3

# Row 7 - new example entry.
$wsMod2.Range("A7").Value = $urlVisibleGroupsA1
$wsMod2.Range("B7").Value = $urlVisibleGroupsA1

# Row 8 - new example entry.
$wsMod2.Range("A8").Value = $urlShow
$wsMod2.Range("B8").Value = $urlShow

# Row 9 - new example entry.
$wsMod2.Range("A9").Value = $urlLat
$wsMod2.Range("B9").Value = $urlLat

# Row 10 - new example entry.
$wsMod2.Range("A10").Value = $urlGclid
$wsMod2.Range("B10").Value = $urlGclid

# Row 11 - mirror of the "Data modifikation" example: short text in A,
# long URL (hyperlinked) in B.
$wsMod2.Range("A11").Value = $shortUrl
$wsMod2.Range("B11").Value = $longUrl
$wsMod2.Hyperlinks.Add($wsMod2.Range("B11"), $longUrl) | Out-Null

$wsMod2.Range("A20").Select() | Out-Null
